# Generate Report for Archive
#
# Two changes:
#   1. Every status cell that read "Ready for handoff" now reads
#      "In Translation" (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
#   2. The Status columns got narrower: Overview columns E & F, and
#      column C on both the "zh-cn" and "de-de" sheets.
#      (The COM layer quantizes ColumnWidth to ~1/6-character steps,
#      so 12.5 is the input that lands closest to the target width.)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn")
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de")

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
